$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A needs the same wrap-text style already used by columns B/E (style index 1)
$ws.Range("A72:A77").WrapText = $true

# Row 72
$ws.Range("A72").Value = "Réalisation"
$ws.Range("B72").Value = "Modification des marges sur la page de login. Les marges posaient des problèmes sur la version mobile de la page car elles étaient calculées par rapport à la taille du viewport"
$ws.Range("C72").Value = 0.5
$ws.Range("D72").Value2 = 43540

# Row 73
$ws.Range("A73").Value = "Analyse"
$ws.Range("B73").Value = "Documentation de l'utilité de mettre en place une standardisation du retour d'erreur de l'API"
$ws.Range("C73").Value = 0.5
$ws.Range("D73").Value2 = 43543

# Row 74
$ws.Range("A74").Value = "Conception"
$ws.Range("B74").Value = "Réflexion sur la structure de retour de l'API ainsi que les types d'erreurs possible"
$ws.Range("C74").Value = 0.5
$ws.Range("D74").Value2 = 43543

# Row 75
$ws.Range("A75").Value = "Réalisation"
$ws.Range("B75").Value = "Création d'une classe de gestion de réponse, modifications des actions utilisateurs pour que les méthodes retournent des Promises"
$ws.Range("C75").Value = 1.5
$ws.Range("D75").Value2 = 43543

# Row 76
$ws.Range("A76").Value = "Réalisation"
$ws.Range("B76").Value = "Transformation des méthodes d'actions sur les activités pour que celles-ci retournent des Promises"
$ws.Range("C76").Value = 0.5
$ws.Range("D76").Value2 = 43543

# Row 77
$ws.Range("A77").Value = "Réalisation"
$ws.Range("B77").Value = "Ajout d'une méthode de login à la classe d'actions utilisalteurs, implémentation de cette méthode lors de l'appel de l'endpoint /token"
$ws.Range("C77").Value = 1
$ws.Range("D77").Value2 = 43543
$ws.Range("E77").Value = "J'ai passé un peu de temps à lire la documentation sur les Promises, je n'était pas totalement au clair avec le rejet des Promises"

# Match the row heights produced by Excel's wrap-text autofit for the new content
$ws.Rows.Item(72).RowHeight = 90
$ws.Rows.Item(73).RowHeight = 60
$ws.Rows.Item(74).RowHeight = 45
$ws.Rows.Item(75).RowHeight = 75
$ws.Rows.Item(76).RowHeight = 60
$ws.Rows.Item(77).RowHeight = 75

# Update frozen pane / view to scroll to the new rows
$ws.Range("A78").Select()

